$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.943378686904907
$ws.Range("B1").Value = 1.957891941070557
$ws.Range("C1").Value = 7.953717708587646
$ws.Range("D1").Value = 0.9170961380004883
$ws.Range("E1").Value = 0.4319026172161102
